$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A89").Value = 3.019
$ws.Range("B89").Value = 1.621
$ws.Range("C89").Value = 2.562
